$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that should no longer hold a value ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()
$ws.Range("B19:C20").Clear()

# --- Row 10: corrected "Objetivos" description ---
$ws.Range("B10:C10").Value = "Apresentar aos estudantes de engenharia os conceitos básicos de Ciência dos Materiais."

# --- Rows 13-15: teacher names (moved under "Docentes responsaveis:") ---
$ws.Range("B13:C13").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("B14:C14").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("B15:C15").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# --- Row 16: "Programa resumido:" + short syllabus text ---
$ws.Range("A16").Value = "Programa resumido:"
$ws.Range("B16:C16").Value = "Estrutura e ligação atômica. 2  Estruturas dos materiais. 3  Imperfeições em sólidos. 4  Diagrama de fases. 5  Propriedades mecânicas"

# --- Row 17: "Short syllabus:" label ---
$ws.Range("A17").Value = "Short syllabus:"

# --- Row 18: "Programa:" + full syllabus text ---
$ws.Range("A18").Value = "Programa:"
$ws.Range("B18:C18").Value = "1. Estrutura e ligação atômica: estrutura dos átomos; ligações covalente, iônica, metálica e forças de van der Waals.`n2. Estruturas dos materiais: sólidos cristalinos; direções e planos cristalográficos; células unitárias; redes de Bravais; fator de empacotamento; métodos para determinação das estruturas cristalinas; estruturas metálicas, iônicas e moleculares. Estrutura de cerâmicas. Estrutura de polímeros. Sólidos amorfos: vidros e polímeros. Aspectos básicos de materiais compósitos. Exemplos de materiais de engenharia.`n3. Imperfeições em sólidos: tipos e formação de defeitos; lacunas; soluções sólidas (intersticial e substitucional); estruturas ordenadas; compostos intermetálicos; discordâncias; movimento de discordâncias; defeitos planares (interfaces). Exemplos práticos.`n4. Diagrama de fases: definição de fase; regra de Gibbs; curva de resfriamento; diagramas de equilíbrio de sistemas binários; equilíbrio de formação e decomposição de fases. Exemplos de diagramas de fases relacionados com a microestrutura dos materiais.`n5. Conceitos básicos sobre as propriedades mecânicas dos materiais: conceitos de tensão e deformação; propriedades elásticas; deformação plástica; plasticidade e fluxo; materiais não newtonianos; relaxação e fluência; fadiga. Exemplos e casos práticos."

# --- Row 19-20: relabeled ---
$ws.Range("A19").Value = "Syllabus:"
$ws.Range("A20").Value = "Avaliação:"

# --- Row 21: "Metodo:" + text ---
$ws.Range("A21").Value = "Método:"
$ws.Range("B21:C21").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."

# --- Row 22: "Criterio:" + text (new row) ---
$ws.Range("A22").Value = "Critério:"
$ws.Range("B22:C22").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."

# --- Row 23: "Norma de recuperacao:" + text (new row) ---
$ws.Range("A23").Value = "Norma de recuperação:"
$ws.Range("B23:C23").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."

# --- Row 24: "Bibliografia:" + text (new row) ---
$ws.Range("A24").Value = "Bibliografia:"
$ws.Range("B24:C24").Value = "1) Askeland, D. R.; Phule, P. P. Ciência e engenharia dos materiais. São Paulo: CENGAGE, 2008.`n2) Callister Jr., W. D. Fundamentos da ciência e engenharia de materiais. Rio de Janeiro: LTC Editora, 2006.`n3) Callister Jr., W. D. Ciência e engenharia de materiais. Rio de Janeiro: LTC Editora, 2008.`n4) Van Vlack, L. H. Princípios de ciência e tecnologia dos materiais. Rio de Janeiro: Editora Campus, 1984.`n5) Shackelford, J. E. Ciência dos materiais. São Paulo: Prentice Hall, 2008. `n6) Jastrzebski, Z. D. The nature and properties of engineering materials. Nova Iorque: John Wiley, 1987.`n7) Padilha, A. F. Materiais de engenharia: microestrutura e propriedades. São Paulo: Hemus Editora, 1997.`n8) Ashby, M. F.; Jones, D. R. H. Engenharia de materiais, 2 vol. Rio de Janeiro: Elsevier Editora, 2007."

# --- Fix formatting for newly created column-B cells: the shared
# "min=1 max=2" / "min=2 max=2" column-B definitions make the engine
# apply column As style (s=1) to brand-new column-B cells instead of
# the correct wrap-text style (s=2). Copy formats from an existing,
# correctly-styled column-B cell to fix this without touching values.
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows("10:11").RowHeight = 60
$ws.Rows("12:15").AutoFit()
$ws.Rows("16:17").RowHeight = 60
$ws.Rows("18:19").RowHeight = 120
$ws.Rows("20:20").AutoFit()
$ws.Rows("21:23").RowHeight = 60
$ws.Rows("24:24").RowHeight = 120

# --- Column A/B split: touching column Bs width causes the engine to
# split the merged "min=1 max=2" column definition, leaving column A
# with its own single-column entry (matching the target diff) while
# column B keeps its already-correct width. ---
$ws.Columns("B").ColumnWidth = $ws.Columns("B").ColumnWidth
